# إضافة حدث جديد في Card9
# 1) Normalize the previously-blank cells on row 25 to contain the literal
#    text "nan" (matching the pattern used by every other data row on this
#    sheet), and
# 2) Append a brand-new service-log entry as row 26.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card9")

# --- Row 25: fill the empty placeholder cells with "nan" ---
$ws.Range("B25:K25").Value = "nan"
$ws.Range("M25").Value = "nan"

# --- Row 26: new service event ---
# Column A repeats the card number ("9"), stored as text like the rest of
# column A, so force text formatting before assigning it, then restore the
# cell style so no visible formatting change is left behind.
$ws.Range("A26").NumberFormat = "@"
$ws.Range("A26").Value = "9"
$ws.Range("A26").NumberFormat = "General"
$ws.Range("A26").Style = "Normal"

$ws.Range("L26").Value = "19\7\2025"
$ws.Range("N26").Value = "تم تغيير الفلاتس وجريده 1"
$ws.Range("O26").Value = "الخبير محمود رشدي"
